# TEST EXECUTION.xlsx - add a new test case row (EXE-06 / EXE-TC-08) to the
# "EXE Regisration" sheet, covering login-after-registration verification.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 is a normal (non-merged) data row that uses exactly the cell styles
# we need for the new row (fills/borders/alignment for columns D:Q). Copy its
# formatting down into row 12 before filling in the new row's values, so the
# new row visually matches the rest of the table without creating any new
# cell-style entries.
$ws.Range("D7:Q7").Copy()
$ws.Range("D12:Q12").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("D12").Value = "EXE-06"
$ws.Range("E12").Value = "Verify user can login after successful registration"
$ws.Range("F12").Value = "EXE-TC-08"
$ws.Range("G12").Value = "Registration Module"
$ws.Range("H12").Value = "EXE-06"
$ws.Range("I12").Value = "Login after successful registration"
$ws.Range("J12").Value = "Login Successful"
$ws.Range("K12").Value = "Login Successful"
$ws.Range("L12").Value = "PASS"
$ws.Range("M12").Value = "-"
$ws.Range("N12").Value = "-"
$ws.Range("O12").Value = "Syaif (QA)"
$ws.Range("P12").Value = 46077
$ws.Range("Q12").Value = "Chrome v145 /`nWindows 18"

# Match the row height used by the other wrapped-text rows in the table.
$ws.Rows.Item(12).RowHeight = 31.5

# Move the active selection down to the newly added row, like the author did.
$ws.Range("K12").Select()
